# Atualização de bases das ligas, do dia: 28-05-2024 às 20:56
#
# This script re-applies a "diff" to the "Northern Ireland Premier" sheet
# that re-shuffles the betting-odds data rows so that each match row
# (everything except the row-index column A, the Div column C and the
# Date column D) ends up holding the data of a different match:
#   - rows 13..17 are cyclically rotated (13->14->15->16->17->13)
#   - rows 137 and 138 are swapped with each other
#   - rows 192 and 193 are swapped with each other
#
# Rather than hard-coding the destination values, we read the current
# ("before") values for every affected cell first and then write them
# back out in the new row order. That way the script is self-contained
# and does not depend on any particular starting state beyond what is
# already in the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that hold data for each match (A/C/D - id/Div/Date - never move)
$cols = @('B','E','F','G','H','I','J','K','L','M','N','O','P','Q','R','S','T','U','V','W','X','Y','Z','AA','AB','AC','AD')

function Get-RowData($rowNum) {
    $data = @{}
    foreach ($c in $cols) {
        $data[$c] = $ws.Range("$c$rowNum").Value()
    }
    return $data
}

function Set-RowData($rowNum, $data) {
    foreach ($c in $cols) {
        $ws.Range("$c$rowNum").Value = $data[$c]
    }
}

# --- Snapshot all the "before" rows we are about to touch ---
$row13 = Get-RowData 13
$row14 = Get-RowData 14
$row15 = Get-RowData 15
$row16 = Get-RowData 16
$row17 = Get-RowData 17

$row137 = Get-RowData 137
$row138 = Get-RowData 138

$row192 = Get-RowData 192
$row193 = Get-RowData 193

# --- Cyclic rotation of rows 13-17: row N receives what used to be in row N-1,
#     and row 13 receives what used to be in row 17 ---
Set-RowData 14 $row13
Set-RowData 15 $row14
Set-RowData 16 $row15
Set-RowData 17 $row16
Set-RowData 13 $row17

# --- Swap rows 137 and 138 ---
Set-RowData 137 $row138
Set-RowData 138 $row137

# --- Swap rows 192 and 193 ---
Set-RowData 192 $row193
Set-RowData 193 $row192
